$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 147
$ws1.Range("F5").Value = 3038
$ws1.Range("F6").Value = 308
$ws1.Range("F7").Value = 407

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 147
$ws4.Range("F5").Value = 3038
$ws4.Range("F6").Value = 308
